$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '82.259.90'
$ws.Range("E2").Value = '  +3.49%  '

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.180.94'
$ws.Range("E3").Value = '  -0.21%  '

# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.00'
$ws.Range("E4").Value = '  +0.07%  '

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '216.32'
$ws.Range("E5").Value = '  +6.09%  '

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '622.83'
$ws.Range("E6").Value = '  -2.04%  '

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.291'
$ws.Range("E7").Value = '  +17.64%  '

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '1.00'
$ws.Range("E8").Value = '  +0.00%  '

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.584'
$ws.Range("E9").Value = '  -3.47%  '

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '3.179.83'
$ws.Range("E10").Value = '  -0.13%  '

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.611'
$ws.Range("E11").Value = '  -3.78%  '

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.0000258'
$ws.Range("E12").Value = '  +5.66%  '

# Row 13
$ws.Range("E13").Value = '  -0.10%  '

# Row 14
$ws.Range("B14").Value = 'WrappedliquidstakedEther2.0'
$ws.Range("C14").Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '3.769.51'
$ws.Range("E14").Value = '  +0.07%  '

# Row 15
$ws.Range("B15").Value = 'Toncoin'
$ws.Range("C15").Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '5.30'
$ws.Range("E15").Value = '  -2.22%  '

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '32.44'
$ws.Range("E16").Value = '  +1.20%  '

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '82.127.41'
$ws.Range("E17").Value = '  +3.62%  '

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '3.181.49'
$ws.Range("E18").Value = '  -0.27%  '

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '3.20'
$ws.Range("E19").Value = '  +7.11%  '

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '14.04'
$ws.Range("E20").Value = '  -3.96%  '

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '435.31'
$ws.Range("E21").Value = '  +0.11%  '

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '8.95'
$ws.Range("E22").Value = '  -3.44%  '

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '5.15'
$ws.Range("E23").Value = '  -0.76%  '

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '7.28'
$ws.Range("E24").Value = '  +5.85%  '

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '5.34'
$ws.Range("E25").Value = '  +9.55%  '

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '11.40'
$ws.Range("E26").Value = '  +3.77%  '

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '3.349.14'
$ws.Range("E27").Value = '  -0.17%  '

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '76.86'
$ws.Range("E28").Value = '  -0.70%  '

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '0.999'
$ws.Range("E29").Value = '  -0.12%  '

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.0000123'
$ws.Range("E30").Value = '  +1.67%  '

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '9.11'
$ws.Range("E31").Value = '  +0.04%  '

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '580.93'
$ws.Range("E32").Value = '  +10.76%  '

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.996'
$ws.Range("E33").Value = '  +0.04%  '

# Row 34
$ws.Range("B34").Value = 'Cronos'
$ws.Range("C34").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.156'
$ws.Range("E34").Value = '  +27.00%  '

# Row 35
$ws.Range("B35").Value = 'Fetch.AI'
$ws.Range("C35").Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.51'
$ws.Range("E35").Value = '  +1.94%  '

# Row 36
$ws.Range("E36").Value = '  +8.52%  '

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '1.99'
$ws.Range("E37").Value = '  +0.03%  '

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '22.83'
$ws.Range("E38").Value = '  -2.17%  '

# Row 39
$ws.Range("B39").Value = 'FirstDigitalUSD'
$ws.Range("C39").Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '1.00'
$ws.Range("E39").Value = '  +0.14%  '

# Row 40
$ws.Range("B40").Value = 'RenderToken'
$ws.Range("C40").Value = 'https://coinranking.com/coin/vfo5XYwcV+rendertoken-render'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '6.16'
$ws.Range("E40").Value = '  +11.58%  '

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.407'
$ws.Range("E41").Value = '  -0.66%  '

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '3.07'
$ws.Range("E42").Value = '  +20.80%  '

# Row 43
$ws.Range("B43").Value = 'WhiteBITCoin'
$ws.Range("C43").Value = 'https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '20.85'
$ws.Range("E43").Value = '  +4.13%  '

# Row 44
$ws.Range("B44").Value = 'Stacks'
$ws.Range("C44").Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '2.02'
$ws.Range("E44").Value = '  +13.06%  '

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '160.94'
$ws.Range("E45").Value = '  -2.55%  '

# Row 46
$ws.Range("E46").Value = '  +0.02%  '

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '188.19'
$ws.Range("E47").Value = '  -1.57%  '

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '44.72'
$ws.Range("E48").Value = '  +2.78%  '

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.33'
$ws.Range("E49").Value = '  +0.67%  '

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '26.53'
$ws.Range("E50").Value = '  +2.97%  '

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.766'
$ws.Range("E51").Value = '  -4.35%  '
